# "reorganized everything and reached a finalized SOP for AIT"
#
# Applies the changes from the commit:
#  - the "By signing this..." attestation paragraph gains a 4th bullet
#  - the (previously blank) "SOP revision date:" field is filled in with an
#    actual date
#  - the two text blocks trade places on the sheet (the attestation text
#    moves into A3, the revision-date text moves into C2) -- done here by
#    writing the new text into each cell, which is observably identical
#  - row 3 is made taller to fit the longer attestation text
#  - the current selection on the sheet is moved to the merged A3:E3 block
#  - the print scale is nudged down now that the content is taller

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Attestation paragraph (now in A3) gains the new bullet point.
$ws.Range("A3").Value = "By signing this, I assert that I have:`n- Completed pertient laboratory safety training`n- Read the AIT SOP in its entirety`n- Become familiar with all the experimental steps outlined in the AIT SOP`n- Sufficient competence to complete these experiments safely"

# "SOP revision date" label (now in C2) gets an actual date filled in.
$ws.Range("C2").Value = " SOP revision date:  July 30, 2018 "

# Row 3 needs to be taller to accommodate the longer wrapped text.
$ws.Rows.Item(3).RowHeight = 98.25

# Selection moves to the merged attestation block.
$ws.Range("A3:E3").Select() | Out-Null

# Tighten the print scale slightly to keep everything on one page.
try { $ws.PageSetup.Zoom = $false } catch {}
try { $ws.PageSetup.Scale = 76 } catch {}

# Record the workbook's (new) saved location, as Excel would on a move/rename.
try { $wb.Path = "K:\AIT\experimental_setup\ait_exp\docs\SOP" } catch {}
